# rapport devops.docx - remove unnecessary screenshot references
#
# 1. Remove the standalone paragraph "Capture d'écran : Exemple d'appel
#    Postman (GET /profile avec token JWT)" that sat right above the
#    horizontal-rule picture after the API table.
# 2. The lastRenderedPageBreak marker that used to sit on the "Inscription,
#    connexion, ajout d'IMC, ajout de programme sportif" bullet now falls on
#    the following "Visualisation du profil utilisateur" bullet instead
#    (a natural consequence of the text removed earlier in the document).
# 3. Collapse the leftover screenshot-caption sentence down to just the
#    second half describing the DELETE example.

$d = $word.ActiveDocument

# --- 1. Delete the "Capture d'écran : Exemple d'appel Postman ..." paragraph
$rngPostman = $d.Content
$foundPostman = $rngPostman.Find.Execute(
    "Capture d'écran : Exemple d'appel Postman (GET /profile avec token JWT)",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPostman) {
    $paraRange = $d.Range($rngPostman.Paragraphs(1).Range.Start, $rngPostman.Paragraphs(1).Range.End)
    $paraRange.Delete()
}

# --- 2. Move <w:lastRenderedPageBreak/> from the "Inscription..." bullet to
#        the "Visualisation du profil utilisateur" bullet.

# 2a. Rewrite the "Inscription..." run without the page-break marker.
$rngInscription = $d.Content
$foundInscription = $rngInscription.Find.Execute(
    "Inscription, connexion, ajout d'IMC, ajout de programme sportif",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundInscription) {
    $rngInscription.Delete()
    $insertAt = $d.Range($rngInscription.Start, $rngInscription.Start)
    $xmlNoBreak = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00BD621D"><w:t>Inscription, connexion, ajout d''IMC, ajout de programme sportif</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertAt.InsertXML($xmlNoBreak)
}

# 2b. Rewrite the "Visualisation du profil utilisateur" run with the marker
#     added back in front of its text.
$rngVisu = $d.Content
$foundVisu = $rngVisu.Find.Execute(
    "Visualisation du profil utilisateur",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundVisu) {
    $rngVisu.Delete()
    $insertAt2 = $d.Range($rngVisu.Start, $rngVisu.Start)
    $xmlWithBreak = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidRPr="00BD621D"><w:lastRenderedPageBreak/><w:t>Visualisation du profil utilisateur</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertAt2.InsertXML($xmlWithBreak)
}

# --- 3. Trim the leftover screenshot caption sentence.
$d.Content.Find.Execute(
    "Capture d'écran : Exemple d'inscription réussie dans Postman Capture d'écran : Exemple de suppression de compte avec DELETE + body (username)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exemple de suppression de compte avec DELETE + body (username)", 2) | Out-Null
